$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (row 1) with two new parallel-line columns: P1=14, Q1=15.
# Copy the existing O1 formatting (bold / centered / top-aligned / bordered, style index 1)
# onto the new cells before assigning their values so the style is reused rather than
# a brand-new style being created.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For every data row (2-25): swap the I/K values, swap the M/O values, and append
# the new P (=2) and Q (=2) columns.
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2
    $oVal = $ws.Cells.Item($r, 15).Value2
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
